# Fruta / hortaliza, semanal
# A new weekly data point is inserted into the time series at row 344,
# pushing the existing rows 344-400 down to 345-401.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 344 (shifts rows 344:400 down to 345:401)
$ws.Rows.Item(344).Insert()

# Populate the newly inserted row 344 with the new record's data.
$ws.Range("A344").Value2 = 9
$ws.Range("B344").Value2 = "Vega Central Mapocho de Santiago"
$ws.Range("C344").Value2 = "Metropolitana"
$ws.Range("D344").Value2 = 44984
$ws.Range("E344").Value2 = 13
$ws.Range("F344").Value2 = 100112043
$ws.Range("G344").Value2 = "Pepino ensalada"
$ws.Range("H344").Value2 = "Sin especificar"
$ws.Range("I344").Value2 = "Primera"
$ws.Range("J344").Value2 = 90
$ws.Range("K344").Value2 = 11000
$ws.Range("L344").Value2 = 13000
$ws.Range("M344").Value2 = 12000
$ws.Range("N344").Value2 = "`$/caja 60 unidades"
$ws.Range("O344").Value2 = "Región de Arica y Parinacota"
$ws.Range("P344").Value2 = 200
$ws.Range("Q344").Value2 = 60
$ws.Range("R344").Value2 = "Hortaliza"
